$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.389662504196167
$ws.Range("B1").Value = 2.901319980621338
$ws.Range("C1").Value = 5.153129577636719
$ws.Range("D1").Value = 1.930641174316406
$ws.Range("E1").Value = 1.205399036407471
